$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GroupPermissionTest")
$c = $ws.Range("A9")
Write-Host $c.Interior.ColorIndex
$c2 = $ws.Range("A2")
Write-Host $c2.Interior.Color
